$d = $word.ActiveDocument

# --- Edit 1 -----------------------------------------------------------
# "Utilized Hyper Markup Language HTML to render web pages using visual
# studio, github with tags that build contents of web pages"
#   -> "Created responsive web pages utilizing Hypertext Markup Language
#       HTML5 and CSS3 using Visual Studio, Github"
$found1 = $d.Content.Find.Execute(
    "Utilized Hyper Markup Language HTML to render web pages using visual studio, github with tags that build contents of web pages",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Created responsive web pages utilizing Hypertext Markup Language HTML5 and CSS3 using Visual Studio, Github",
    2)
Write-Host "Edit1 replaced: $found1"

# --- Edit 2 -------------------------------------------------------------
# The second bullet's paragraph is split across two runs: a lone leading
# "U" run (carrying extra w:color / w:vertAlign formatting) immediately
# followed by a run starting with "Utilized Hyper Markup Language HTML to
# render web pages with tags that build contents of web pages" -- so the
# paragraph literally reads "UUtilized...". The target keeps only the
# second run (trimmed rPr, no color/vertAlign) with updated wording, and
# drops the stray leading "U" run entirely. We find the paragraph
# dynamically (rather than relying on fixed character offsets, which
# shift once Edit 1 has run) and delete just its first character.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -like "UUtilized Hyper Markup Language HTML to render web pages with tags*") {
        $startPos = $p.Range.Start
        $leadChar = $d.Range($startPos, $startPos + 1)
        $leadChar.Delete()
        break
    }
}

# Now the paragraph reads "Utilized Hyper Markup Language HTML to render
# web pages with tags that build contents of web pages" in a single run;
# update its wording in place.
$found2 = $d.Content.Find.Execute(
    "Utilized Hyper Markup Language HTML to render web pages with tags that build contents of web pages",
    $true, $false, $false, $false, $false, $true, 1, $false,
    "Utilized Hyper Markup Language HTML to create web pages with tags",
    2)
Write-Host "Edit2 replaced: $found2"
